# Auto-generated script applying Phantom_Profits.xlsx data refresh
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(9, 8).Value = 945.9091
$ws.Cells.Item(9, 9).Value = 886.5714
$ws.Cells.Item(9, 10).Value = 1049.75
$ws.Cells.Item(9, 11).Value = 886.5714
$ws.Cells.Item(9, 12).Value = 1049.75
$ws.Cells.Item(9, 13).Value = -717.5714
$ws.Cells.Item(9, 14).Value = -1387.75

$ws.Cells.Item(17, 8).Value = 1695.4333
$ws.Cells.Item(17, 10).Value = 1695.4333
$ws.Cells.Item(17, 12).Value = 5086.2999
$ws.Cells.Item(17, 14).Value = -5422.2999

$ws.Cells.Item(28, 8).Value = 1777.7693
$ws.Cells.Item(28, 9).Value = 1160.5
$ws.Cells.Item(28, 10).Value = 3835.3333
$ws.Cells.Item(28, 11).Value = 1160.5
$ws.Cells.Item(28, 12).Value = 3835.3333
$ws.Cells.Item(28, 13).Value = -675.5
$ws.Cells.Item(28, 14).Value = -4805.3333

$ws.Cells.Item(29, 8).Value = 1308.5714
$ws.Cells.Item(29, 9).Value = 1110
$ws.Cells.Item(29, 11).Value = 3330
$ws.Cells.Item(29, 13).Value = -3049

$ws.Cells.Item(38, 8).Value = 1250221.9
$ws.Cells.Item(38, 10).Value = 0
$ws.Cells.Item(38, 12).Value = 0
$ws.Cells.Item(38, 14).ClearContents()

$ws.Cells.Item(58, 8).Value = 527.5
$ws.Cells.Item(58, 9).Value = 527.5
$ws.Cells.Item(58, 10).Value = 0
$ws.Cells.Item(58, 11).Value = 1582.5
$ws.Cells.Item(58, 12).Value = 0
$ws.Cells.Item(58, 13).Value = -1432.5
$ws.Cells.Item(58, 14).ClearContents()

$ws.Cells.Item(74, 8).Value = 3603.6365
$ws.Cells.Item(74, 9).Value = 3603.6365
$ws.Cells.Item(74, 11).Value = 3603.6365
$ws.Cells.Item(74, 13).Value = -2667.6365

$ws.Cells.Item(77, 8).Value = 3603.6365
$ws.Cells.Item(77, 9).Value = 3603.6365
$ws.Cells.Item(77, 11).Value = 18018.1825
$ws.Cells.Item(77, 13).Value = -13338.1825

$ws.Cells.Item(98, 8).Value = 854.1
$ws.Cells.Item(98, 9).Value = 854.1
$ws.Cells.Item(98, 11).Value = 854.1
$ws.Cells.Item(98, 13).Value = 643.9

$ws.Cells.Item(100, 8).Value = 2337.5881
$ws.Cells.Item(100, 9).Value = 2189.3333
$ws.Cells.Item(100, 11).Value = 2189.3333
$ws.Cells.Item(100, 13).Value = -1648.3333

$ws.Cells.Item(122, 8).Value = 854.1
$ws.Cells.Item(122, 9).Value = 854.1
$ws.Cells.Item(122, 11).Value = 2562.3
$ws.Cells.Item(122, 13).Value = -112.3000000000002

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 2863.9285
$ws.Cells.Item(61, 9).Value = 2868.8462
$ws.Cells.Item(61, 11).Value = 2868.8462
$ws.Cells.Item(61, 13).Value = -2656.8462

$ws.Cells.Item(97, 8).Value = 499.33334
$ws.Cells.Item(97, 9).Value = 499.33334
$ws.Cells.Item(97, 11).Value = 499.33334
$ws.Cells.Item(97, 13).Value = -3.333340000000021

$ws.Cells.Item(136, 8).Value = 2863.9285
$ws.Cells.Item(136, 9).Value = 2868.8462
$ws.Cells.Item(136, 11).Value = 8606.5386
$ws.Cells.Item(136, 13).Value = -6056.5386

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(11, 8).Value = 335
$ws.Cells.Item(11, 9).Value = 4
$ws.Cells.Item(11, 10).Value = 500.5
$ws.Cells.Item(11, 11).Value = 4
$ws.Cells.Item(11, 12).Value = 500.5
$ws.Cells.Item(11, 13).Value = 136
$ws.Cells.Item(11, 14).Value = -780.5

$ws.Cells.Item(94, 8).Value = 428.47058
$ws.Cells.Item(94, 10).Value = 541.8570999999999
$ws.Cells.Item(94, 12).Value = 541.8570999999999
$ws.Cells.Item(94, 14).Value = -1443.8571

$ws.Cells.Item(107, 8).Value = 828.61536
$ws.Cells.Item(107, 9).Value = 647.75
$ws.Cells.Item(107, 10).Value = 2999
$ws.Cells.Item(107, 11).Value = 647.75
$ws.Cells.Item(107, 12).Value = 2999
$ws.Cells.Item(107, 13).Value = 1272.25
$ws.Cells.Item(107, 14).Value = -6839

$ws.Cells.Item(134, 8).Value = 3993.5789
$ws.Cells.Item(134, 9).Value = 4021
$ws.Cells.Item(134, 10).Value = 3500
$ws.Cells.Item(134, 11).Value = 12063
$ws.Cells.Item(134, 12).Value = 10500
$ws.Cells.Item(134, 13).Value = -9528
$ws.Cells.Item(134, 14).Value = -15570

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 386.1111
$ws.Cells.Item(7, 9).Value = 308
$ws.Cells.Item(7, 11).Value = 308
$ws.Cells.Item(7, 13).Value = -195

$ws.Cells.Item(99, 8).Value = 3021.3333
$ws.Cells.Item(99, 9).Value = 2899.6667
$ws.Cells.Item(99, 11).Value = 2899.6667
$ws.Cells.Item(99, 13).Value = -1401.6667

$ws.Cells.Item(107, 8).Value = 1630
$ws.Cells.Item(107, 10).Value = 1630
$ws.Cells.Item(107, 12).Value = 1630
$ws.Cells.Item(107, 14).Value = -5470

$ws.Cells.Item(126, 8).Value = 3021.3333
$ws.Cells.Item(126, 9).Value = 2899.6667
$ws.Cells.Item(126, 11).Value = 8699.000100000001
$ws.Cells.Item(126, 13).Value = -6229.000100000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(75, 8).Value = 4062
$ws.Cells.Item(75, 10).Value = 2932.3333
$ws.Cells.Item(75, 12).Value = 8796.999899999999
$ws.Cells.Item(75, 14).Value = -10792.9999

$ws.Cells.Item(78, 8).Value = 4062
$ws.Cells.Item(78, 10).Value = 2932.3333
$ws.Cells.Item(78, 12).Value = 26390.9997
$ws.Cells.Item(78, 14).Value = -36374.9997

$ws.Cells.Item(86, 8).Value = 1065.875
$ws.Cells.Item(86, 10).Value = 1994.3334
$ws.Cells.Item(86, 12).Value = 5983.0002
$ws.Cells.Item(86, 14).Value = -8355.0002

$ws.Cells.Item(89, 8).Value = 1065.875
$ws.Cells.Item(89, 10).Value = 1994.3334
$ws.Cells.Item(89, 12).Value = 17949.0006
$ws.Cells.Item(89, 14).Value = -29805.0006

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(126, 8).Value = 1816.5
$ws.Cells.Item(126, 9).Value = 1816.5
$ws.Cells.Item(126, 10).Value = 0
$ws.Cells.Item(126, 11).Value = 5449.5
$ws.Cells.Item(126, 12).Value = 0
$ws.Cells.Item(126, 13).Value = -2979.5
$ws.Cells.Item(126, 14).ClearContents()

$ws.Cells.Item(132, 8).Value = 2814.7896
$ws.Cells.Item(132, 9).Value = 2872.75
$ws.Cells.Item(132, 10).Value = 2772.6365
$ws.Cells.Item(132, 11).Value = 8618.25
$ws.Cells.Item(132, 12).Value = 8317.9095
$ws.Cells.Item(132, 13).Value = -6088.25
$ws.Cells.Item(132, 14).Value = -13377.9095

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 4949.3687
$ws.Cells.Item(7, 9).Value = 5007.6665
$ws.Cells.Item(7, 10).Value = 3900
$ws.Cells.Item(7, 11).Value = 5007.6665
$ws.Cells.Item(7, 12).Value = 3900
$ws.Cells.Item(7, 13).Value = -4895.6665
$ws.Cells.Item(7, 14).Value = -4124

$ws.Cells.Item(40, 8).Value = 3973.75
$ws.Cells.Item(40, 9).Value = 3973.75
$ws.Cells.Item(40, 10).Value = 0
$ws.Cells.Item(40, 11).Value = 3973.75
$ws.Cells.Item(40, 12).Value = 0
$ws.Cells.Item(40, 13).Value = -3837.75
$ws.Cells.Item(40, 14).ClearContents()

$ws.Cells.Item(46, 8).Value = 1370.4445
$ws.Cells.Item(46, 9).Value = 1330
$ws.Cells.Item(46, 10).Value = 1421
$ws.Cells.Item(46, 11).Value = 1330
$ws.Cells.Item(46, 12).Value = 1421
$ws.Cells.Item(46, 13).Value = -1142
$ws.Cells.Item(46, 14).Value = -1797

$ws.Cells.Item(55, 8).Value = 367.25
$ws.Cells.Item(55, 9).Value = 305.66666
$ws.Cells.Item(55, 10).Value = 446.42856
$ws.Cells.Item(55, 11).Value = 305.66666
$ws.Cells.Item(55, 12).Value = 446.42856
$ws.Cells.Item(55, 13).Value = -132.66666
$ws.Cells.Item(55, 14).Value = -792.4285600000001

$ws.Cells.Item(122, 8).Value = 3548.35
$ws.Cells.Item(122, 9).Value = 3204.25
$ws.Cells.Item(122, 10).Value = 4924.75
$ws.Cells.Item(122, 11).Value = 9612.75
$ws.Cells.Item(122, 12).Value = 14774.25
$ws.Cells.Item(122, 13).Value = -7162.75
$ws.Cells.Item(122, 14).Value = -19674.25

$ws.Cells.Item(126, 8).Value = 4949.3687
$ws.Cells.Item(126, 9).Value = 5007.6665
$ws.Cells.Item(126, 10).Value = 3900
$ws.Cells.Item(126, 11).Value = 15022.9995
$ws.Cells.Item(126, 12).Value = 11700
$ws.Cells.Item(126, 13).Value = -12552.9995
$ws.Cells.Item(126, 14).Value = -16640

$ws.Cells.Item(132, 8).Value = 5749
$ws.Cells.Item(132, 9).Value = 5749
$ws.Cells.Item(132, 11).Value = 17247
$ws.Cells.Item(132, 13).Value = -14717

$ws.Cells.Item(140, 8).Value = 0
$ws.Cells.Item(140, 10).Value = 0
$ws.Cells.Item(140, 12).Value = 0
$ws.Cells.Item(140, 14).ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(44, 8).Value = 0
$ws.Cells.Item(44, 10).Value = 0
$ws.Cells.Item(44, 12).Value = 0
$ws.Cells.Item(44, 14).ClearContents()

$ws.Cells.Item(113, 8).Value = 207.33333
$ws.Cells.Item(113, 9).Value = 255.4
$ws.Cells.Item(113, 11).Value = 766.2
$ws.Cells.Item(113, 13).Value = 1403.8

$ws.Cells.Item(122, 8).Value = 3057.8518
$ws.Cells.Item(122, 9).Value = 2998.524
$ws.Cells.Item(122, 11).Value = 8995.572
$ws.Cells.Item(122, 13).Value = -6545.572

$ws.Cells.Item(126, 8).Value = 4698.375
$ws.Cells.Item(126, 9).Value = 4670
$ws.Cells.Item(126, 10).Value = 4897
$ws.Cells.Item(126, 11).Value = 14010
$ws.Cells.Item(126, 12).Value = 14691
$ws.Cells.Item(126, 13).Value = -11540
$ws.Cells.Item(126, 14).Value = -19631

$ws.Cells.Item(132, 8).Value = 4670.68
$ws.Cells.Item(132, 9).Value = 4763.35
$ws.Cells.Item(132, 10).Value = 4300
$ws.Cells.Item(132, 11).Value = 14290.05
$ws.Cells.Item(132, 12).Value = 12900
$ws.Cells.Item(132, 13).Value = -11760.05
$ws.Cells.Item(132, 14).Value = -17960

$ws.Cells.Item(136, 8).Value = 6484.914
$ws.Cells.Item(136, 9).Value = 4026.2222
$ws.Cells.Item(136, 10).Value = 9088.235000000001
$ws.Cells.Item(136, 11).Value = 12078.6666
$ws.Cells.Item(136, 12).Value = 27264.705
$ws.Cells.Item(136, 13).Value = -9528.6666
$ws.Cells.Item(136, 14).Value = -32364.705
